$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 422, shifting all following rows down by one.
$ws.Rows.Item(422).Insert()

# Populate the newly inserted row 422 with the new data record.
$ws.Cells.Item(422, 1).Value = 9
$ws.Cells.Item(422, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(422, 3).Value = "Metropolitana"
$ws.Cells.Item(422, 4).Value = 44711
$ws.Cells.Item(422, 5).Value = 13
$ws.Cells.Item(422, 6).Value = 100112031
$ws.Cells.Item(422, 7).Value = "Poroto verde"
$ws.Cells.Item(422, 8).Value = "Magnum"
$ws.Cells.Item(422, 9).Value = "Primera"
$ws.Cells.Item(422, 10).Value = 43
$ws.Cells.Item(422, 11).Value = 28000
$ws.Cells.Item(422, 12).Value = 30000
$ws.Cells.Item(422, 13).Value = 29023
$ws.Cells.Item(422, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(422, 15).Value = "Perú"
$ws.Cells.Item(422, 16).Value = 1161
$ws.Cells.Item(422, 17).Value = 25
$ws.Cells.Item(422, 18).Value = "Hortaliza"
